$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''68.721.14'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '''3.490.07'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''588.28'
$ws.Range('E5').Value = '  +2.26%  '
$ws.Range('D6').Value = '''168.82'
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  -1.48%  '
$ws.Range('D8').Value = '''3.482.80'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').Value = '''1.00'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('D11').Value = '''6.73'
$ws.Range('E11').Value = '  +3.00%  '
$ws.Range('D12').Value = '''0.574'
$ws.Range('E12').Value = '  -3.95%  '
$ws.Range('D13').Value = '''46.76'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').Value = '''4.055.02'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('D16').Value = '''616.71'
$ws.Range('E16').Value = '  -10.28%  '
$ws.Range('D17').Value = '''8.33'
$ws.Range('E17').Value = '  -4.70%  '
$ws.Range('D18').Value = '''3.485.69'
$ws.Range('E18').Value = '  -1.25%  '
$ws.Range('D19').Value = '''68.897.43'
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('E20').Value = '  -2.22%  '
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('E23').Value = '  -3.79%  '
$ws.Range('D24').Value = '''15.83'
$ws.Range('E24').Value = '  -4.05%  '
$ws.Range('D25').Value = '''95.86'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('D26').Value = '''3.78'
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').Value = '''2.61'
$ws.Range('E28').Value = '  -2.51%  '
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('D30').Value = '''32.95'
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('D31').Value = '''8.43'
$ws.Range('E31').Value = '  -4.41%  '
$ws.Range('E32').Value = '  -2.10%  '
$ws.Range('D33').Value = '''1.32'
$ws.Range('E33').Value = '  -2.47%  '
$ws.Range('D34').Value = '''6.85'
$ws.Range('E34').Value = '  -5.38%  '
$ws.Range('D35').Value = '''577.31'
$ws.Range('E35').Value = '  +1.60%  '
$ws.Range('B36').Value = 'dogwifhat'
$ws.Range('C36').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D36').Value = '''3.53'
$ws.Range('E36').Value = '  -5.75%  '
$ws.Range('B37').Value = 'Cosmos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D37').Value = '''10.70'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('D38').Value = '''56.99'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('E39').Value = '  -3.20%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.136'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0435'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').Value = '''3.412.32'
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('E44').Value = '  -3.75%  '
$ws.Range('E45').Value = '  -1.78%  '
$ws.Range('D46').Value = '''0.0₃0693'
$ws.Range('E46').Value = '  -1.40%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').Value = '''2.80'
$ws.Range('E47').Value = '  -3.43%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').Value = '''2.54'
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('E49').Value = '  -2.60%  '
$ws.Range('D50').Value = '''5.74'
$ws.Range('E50').Value = '  +12.87%  '
$ws.Range('D51').Value = '''132.15'
$ws.Range('E51').Value = '  -1.10%  '
